# Commit: [ADDITIONAL SCRAPING] added scraping code for extra browling attributes and excel sheets
#
# 1) "ODI Batting Extra" had a bunch of cells that were present but holding no
#    value (leftover empty inline-string placeholders from the scrape). Those
#    are removed now that the scraper only emits cells it actually has data
#    for.
# 2) A brand-new sheet "ODI Bowling Extra" is added (after "ODI Batting
#    Extra") holding MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL for
#    the bowling-extras scrape.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: drop the now-empty placeholder cells on "ODI Batting Extra"
# ---------------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")

$emptyCells = @(
    "C3","D3","E3",
    "B5","C5","D5","E5",
    "B7","C7","D7","E7",
    "B9","C9","D9","E9",
    "C10","D10","E10",
    "B11","C11","D11","E11",
    "B13","C13","D13","E13",
    "B17","C17","D17","E17",
    "B19","C19","D19","E19",
    "B20","C20","D20","E20"
)

foreach ($addr in $emptyCells) {
    $battingExtra.Range($addr).ClearContents()
}

# ---------------------------------------------------------------------------
# Step 2: add the new "ODI Bowling Extra" sheet as the last tab
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bowlingExtra = $wb.Worksheets.Add($null, $lastSheet)
$bowlingExtra.Name = "ODI Bowling Extra"

# Header row
$bowlingExtra.Range("A1").Value = "MATCH_CODE"
$bowlingExtra.Range("B1").Value = "MAIDEN_OVERS"
$bowlingExtra.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# Match the look of the other "Extra" sheets' header (bold, centered, thin
# border) by copying the formatting from the existing header.
$battingExtra.Range("A1:C1").Copy()
$bowlingExtra.Range("A1:C1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Step 3: populate the scraped rows
# ---------------------------------------------------------------------------
# Columns: MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL
# $null entries mean the scraper had no value for that field for that match.
$rows = @(
    @("4442", "0",    "10.00%"),
    @("4444", "0",    "50.00%"),
    @("4446", "0",    $null),
    @("4448", "2",    "10.00%"),
    @("4466", "1",    "10.00%"),
    @("4467", "4",    "10.00%"),
    @("4468", "1",    "10.00%"),
    @("4475", "0",    "20.00%"),
    @("4478", $null,  $null),
    @("4492", "0",    "10.00%"),
    @("4496", "2",    "30.00%"),
    @("4519", $null,  $null),
    @("4520", "2",    "40.00%"),
    @("4522", $null,  $null),
    @("4605", "0",    "10.00%"),
    @("4608", "0",    $null),
    @("4614", "0",    $null),
    @("4694", "0",    "10.00%"),
    @("4726", $null,  $null),
    @("4729", $null,  $null)
)

# Force every data cell in this block to be stored as text (matching the
# scraper output, which never emits numeric cells on this sheet) instead of
# letting COM auto-detect numbers/percentages.
$dataRange = $bowlingExtra.Range("A2:C" + (1 + $rows.Length))
$dataRange.NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $rowVals = $rows[$i]

    if ($rowVals[0] -ne $null) {
        $bowlingExtra.Range("A$r").Value = $rowVals[0]
    }
    if ($rowVals[1] -ne $null) {
        $bowlingExtra.Range("B$r").Value = $rowVals[1]
    }
    if ($rowVals[2] -ne $null) {
        $bowlingExtra.Range("C$r").Value = $rowVals[2]
    }
}
